$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.528.59"
$ws.Range("E2").Value = "  +5.43%  "

$ws.Range("D3").Value = "2.298.82"
$ws.Range("E3").Value = "  +4.39%  "

$ws.Range("E4").Value = "  -0.76%  "

$ws.Range("D5").Value = "'299.91"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").Value = "'96.42"
$ws.Range("E6").Value = "  +8.39%  "

$ws.Range("D7").Value = "'0.571"
$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("E8").Value = "  -0.57%  "

$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  +7.22%  "

$ws.Range("D10").Value = "'35.51"
$ws.Range("E10").Value = "  +5.58%  "

$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  +2.02%  "

$ws.Range("D12").Value = "'7.40"
$ws.Range("E12").Value = "  +7.26%  "

$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "2.651.20"
$ws.Range("E14").Value = "  +4.25%  "

$ws.Range("D15").Value = "2.302.56"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").Value = "'0.829"
$ws.Range("E16").Value = "  +5.23%  "

$ws.Range("D17").Value = "'13.98"
$ws.Range("E17").Value = "  +5.86%  "

$ws.Range("D18").Value = "46.417.37"
$ws.Range("E18").Value = "  +5.26%  "

$ws.Range("D19").Value = "'13.27"
$ws.Range("E19").Value = "  +18.09%  "

$ws.Range("D20").Value = "0.0₃0944"
$ws.Range("E20").Value = "  +5.29%  "

$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  +3.96%  "

$ws.Range("D22").Value = "'67.16"
$ws.Range("E22").Value = "  +4.65%  "

$ws.Range("D23").Value = "'250.49"
$ws.Range("E23").Value = "  +7.76%  "

$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  +2.82%  "

$ws.Range("D25").Value = "'2.00"
$ws.Range("E25").Value = "  +5.24%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").Value = "'41.73"
$ws.Range("E27").Value = "  +13.51%  "

$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").Value = "'9.83"
$ws.Range("E29").Value = "  +5.13%  "

$ws.Range("D30").Value = "'20.13"
$ws.Range("E30").Value = "  +4.08%  "

$ws.Range("D31").Value = "'5.81"
$ws.Range("E31").Value = "  +5.04%  "

$ws.Range("D32").Value = "'0.0804"
$ws.Range("E32").Value = "  +7.16%  "

$ws.Range("D33").Value = "'145.19"
$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").Value = "'2.60"
$ws.Range("E34").Value = "  +2.78%  "

$ws.Range("D35").Value = "'3.10"
$ws.Range("E35").Value = "  +6.99%  "

$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  +5.62%  "

$ws.Range("E37").Value = "  +1.05%  "

$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  +6.35%  "

$ws.Range("E39").Value = "  +13.07%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'14.75"
$ws.Range("E40").Value = "  +6.27%  "

$ws.Range("D41").Value = "'3.39"
$ws.Range("E41").Value = "  +6.42%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0307"
$ws.Range("E42").Value = "  +8.06%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.83%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.97"
$ws.Range("E44").Value = "  +20.26%  "

$ws.Range("D45").Value = "1.817.96"
$ws.Range("E45").Value = "  +4.03%  "

$ws.Range("D46").Value = "'91.81"
$ws.Range("E46").Value = "  +21.36%  "

$ws.Range("D47").Value = "'0.193"
$ws.Range("E47").Value = "  +9.48%  "

$ws.Range("D48").Value = "'72.75"
$ws.Range("E48").Value = "  +6.83%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'4.86"
$ws.Range("E49").Value = "  +9.31%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'97.36"
$ws.Range("E50").Value = "  +4.31%  "

$ws.Range("D51").Value = "'54.53"
$ws.Range("E51").Value = "  +6.78%  "

